$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 6 (pushes existing rows 6..125 down to 7..126)
$ws.Rows("6:6").Insert()

# Populate the new row 6 with the new price record
$ws.Range("A6").Value = 10
$ws.Range("B6").Value = "Vega Modelo de Temuco"
$ws.Range("C6").Value = "La Araucanía"
$ws.Range("D6").Value = 44817
$ws.Range("E6").Value = 9
$ws.Range("F6").Value = 100114002
$ws.Range("G6").Value = "Camote"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 25
$ws.Range("K6").Value = 20000
$ws.Range("L6").Value = 20000
$ws.Range("M6").Value = 20000
$ws.Range("N6").Value = "$/malla 20 kilos"
$ws.Range("O6").Value = "Perú"
$ws.Range("P6").Value = 1000
$ws.Range("Q6").Value = 20
$ws.Range("R6").Value = "Hortaliza"
